$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $c = $sheet.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" '42.350.66'
Set-TextValue $ws "E2" '  +1.28%  '
Set-TextValue $ws "D3" '2.271.92'
Set-TextValue $ws "E3" '  +0.02%  '
Set-TextValue $ws "E4" '  +0.04%  '
Set-TextValue $ws "D6" '97.35'
Set-TextValue $ws "E6" '  +5.06%  '
Set-TextValue $ws "E7" '  -0.09%  '
Set-TextValue $ws "E8" '  +0.02%  '
Set-TextValue $ws "D9" '0.491'
Set-TextValue $ws "E9" '  +1.27%  '
Set-TextValue $ws "D10" '35.70'
Set-TextValue $ws "E10" '  +9.58%  '
Set-TextValue $ws "D11" '0.0796'
Set-TextValue $ws "E11" '  -0.11%  '
Set-TextValue $ws "E12" '  -0.89%  '
Set-TextValue $ws "D13" '6.66'
Set-TextValue $ws "E13" '  -0.50%  '
Set-TextValue $ws "D14" '2.624.76'
Set-TextValue $ws "E14" '  +0.04%  '
Set-TextValue $ws "E15" '  +0.50%  '
Set-TextValue $ws "D16" '2.278.49'
Set-TextValue $ws "E16" '  +0.97%  '
Set-TextValue $ws "D17" '0.794'
Set-TextValue $ws "E17" '  +1.90%  '
Set-TextValue $ws "D18" '42.244.47'
Set-TextValue $ws "E18" '  +1.24%  '
Set-TextValue $ws "D19" '12.52'
Set-TextValue $ws "E19" '  -0.61%  '
Set-TextValue $ws "D20" '0.0₃0909'
Set-TextValue $ws "E20" '  +0.14%  '
Set-TextValue $ws "E21" '  +0.18%  '
Set-TextValue $ws "E22" '  +0.68%  '
Set-TextValue $ws "D23" '240.39'
Set-TextValue $ws "E23" '  -1.45%  '
Set-TextValue $ws "E24" '  +0.17%  '
Set-TextValue $ws "D25" '1.96'
Set-TextValue $ws "E25" '  +1.22%  '
Set-TextValue $ws "E26" '  -0.17%  '
Set-TextValue $ws "D27" '23.81'
Set-TextValue $ws "E27" '  -0.92%  '
Set-TextValue $ws "D28" '37.26'
Set-TextValue $ws "E28" '  +5.25%  '
Set-TextValue $ws "D29" '9.51'
Set-TextValue $ws "E29" '  +0.01%  '
Set-TextValue $ws "D30" '2.11'
Set-TextValue $ws "E30" '  +1.56%  '
Set-TextValue $ws "D31" '159.82'
Set-TextValue $ws "E31" '  -0.49%  '
Set-TextValue $ws "E32" '  +0.01%  '
Set-TextValue $ws "E33" '  +0.04%  '
Set-TextValue $ws "E34" '  +4.77%  '
Set-TextValue $ws "D35" '0.0742'
Set-TextValue $ws "E35" '  -0.10%  '
Set-TextValue $ws "D36" '17.06'
Set-TextValue $ws "E36" '  +0.47%  '
Set-TextValue $ws "E37" '  +0.37%  '
Set-TextValue $ws "D38" '2.33'
Set-TextValue $ws "E38" '  -1.23%  '
Set-TextValue $ws "E39" '  +1.74%  '
Set-TextValue $ws "E40" '  -1.49%  '
Set-TextValue $ws "D41" '4.07'
Set-TextValue $ws "E41" '  +4.19%  '
Set-TextValue $ws "E42" '  +14.17%  '
Set-TextValue $ws "D43" '1.997.43'
Set-TextValue $ws "E43" '  -0.36%  '
Set-TextValue $ws "E44" '  +1.07%  '
Set-TextValue $ws "D45" '18.91'
Set-TextValue $ws "E45" '  -4.14%  '
Set-TextValue $ws "D46" '2.94'
Set-TextValue $ws "E46" '  +1.40%  '
Set-TextValue $ws "D47" '9.97'
Set-TextValue $ws "E47" '  -3.05%  '
Set-TextValue $ws "D48" '53.22'
Set-TextValue $ws "E48" '  +0.92%  '
Set-TextValue $ws "E49" '  +0.59%  '
Set-TextValue $ws "E50" '  -0.07%  '
Set-TextValue $ws "D51" '91.50'
Set-TextValue $ws "E51" '  +0.39%  '
